# Financial Statement Output.xlsx edit
# Replaces the "balance sheet / income statement" style block (rows 1-95) with
# a much shorter "comprehensive income" reconciliation block (rows 1-9), and
# blanks out the remaining rows 10-95 so they match the already-empty rows
# that follow (96+).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: column headers (period end dates) ---
$ws.Range("D1").Value = "september 29,"
$ws.Range("E1").Value = "september 30, 2017"
$ws.Range("F1").Value = "september 24,"

# --- Row 2 ---
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = " net income"
$ws.Range("D2").Value = 59531
$ws.Range("E2").Value = 48351
$ws.Range("F2").Value = 45687

# --- Row 3 ---
$ws.Range("C3").Value = " respectively"
$ws.Range("D3").Value = -525
$ws.Range("E3").Value = 224
$ws.Range("F3").Value = 75

# --- Row 4 (F4 stays blank, as in the source diff) ---
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = " (478) and (7) respectively"
$ws.Range("D4").Value = 523
$ws.Range("E4").Value = 1315

# --- Row 5 ---
$ws.Range("B5").Value = 11
$ws.Range("C5").Value = " of tax"
$ws.Range("D5").Value = 905
$ws.Range("E5").Value = -162
$ws.Range("F5").Value = -734

# --- Row 6 ---
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = " and (863) respectively"
$ws.Range("D6").Value = -3407
$ws.Range("E6").Value = -782
$ws.Range("F6").Value = 1582

# --- Row 7 ---
$ws.Range("B7").Value = 18
$ws.Range("C7").Value = " tax"
$ws.Range("D7").Value = -3406
$ws.Range("E7").Value = -846
$ws.Range("F7").Value = 1638

# --- Row 8 ---
$ws.Range("B8").Value = 19
$ws.Range("C8").Value = " total other cuniprehensive incuine/(loss)"
$ws.Range("D8").Value = -3026
$ws.Range("E8").Value = -784
$ws.Range("F8").Value = 979

# --- Row 9 ---
$ws.Range("B9").Value = 20
$ws.Range("C9").Value = " total comprehensive income see accompanying apple notes inc. to | consolidated fun 10-k financial | statements"
$ws.Range("D9").Value = 56505
$ws.Range("E9").Value = 47567
$ws.Range("F9").Value = 46666

# --- Rows 10-95: drop the page_num/line_num/variable columns entirely, and
#     blank out the numeric columns (same shape as the already-empty rows
#     96 onward). ---
$ws.Range("A10:F95").ClearContents()
